$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("crossBar")

# PB7's "Used" note becomes more specific now that PE15 also carries a
# timer-related signal: PB7 drives the timer interrupt itself.
$ws.Range("D18").Value = "TIMER DRIVER INT"

# Reserve PE15 (row 29) for the timer integration test output. The cell
# previously only carried a direct border style with no content; clear
# that leftover formatting before writing the new label so it reverts to
# the sheet's default (unstyled) cell format.
$ws.Range("D29").ClearFormats()
$ws.Range("D29").Value = "TIMER INT"

# The "Used" column now needs to fit the longer "TIMER DRIVER INT" label.
$ws.Columns.Item(4).ColumnWidth = 16

# Scroll the saved view down a bit and leave the newly-edited cell selected.
$ws.Range("D29").Select()
$excel.ActiveWindow.ScrollRow = 11
